$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update C2: value changes from 1.2 to 3.2 and loses its custom number-format style
$ws.Range("C2").Value = 3.2
$ws.Range("C2").ClearFormats()

# Append new ranking rows 439-476 (columns A: name, B: points)
$data = @(
  @(439, "רומי הרשקוביץ", 1),
  @(440, "איתי הראל", 1),
  @(441, "ליהי בראל", 1),
  @(442, "אורי שטרנברג", 1),
  @(443, "ירון גלפנד", 1),
  @(444, "אורי שטרנברג", 6),
  @(445, "ירון גלפנד", 6),
  @(446, "עדן ורד מרי", 1),
  @(447, "אביב ואסקז", 1),
  @(448, "הילה שולויס", 1),
  @(449, "יהלי גודר", 1),
  @(450, "ליאם דיין", 1),
  @(451, "תומר ששון", 1),
  @(452, "אן מרש", 1),
  @(453, "יולי קזמה", 1),
  @(454, "איתי בסטקר", 1),
  @(455, "מעיין סטרוזר", 1),
  @(456, "יהלי דוייב", 1),
  @(457, "הילה שולויס", 6),
  @(458, "עדן ורד מרי", 6),
  @(459, "אביב ואסקז", 1),
  @(460, "יהלי דוייב", 1),
  @(461, "תומר ששון", 1),
  @(462, "תאיו ורד", 1),
  @(463, "שלו דיין", 1),
  @(464, "אביב ואסקז", 6),
  @(465, "יהלי דוייב", 6),
  @(466, "רומי הרשקוביץ", 1),
  @(467, "עדן ורד מרי", 1),
  @(468, "מעיין סטרוזר", 1),
  @(469, "הילה שולויס", 1),
  @(470, "ליאם דיין", 1),
  @(471, "יער אלביר", 1),
  @(472, "יולי קזמה", 1),
  @(473, "ליהי בראל", 1),
  @(474, "אורי שטרנברג", 1),
  @(475, "ליאם דיין", 6),
  @(476, "רומי הרשקוביץ", 6)
)

foreach ($row in $data) {
    $r = $row[0]
    $name = $row[1]
    $val = $row[2]
    $ws.Cells.Item($r, 1).Value = $name
    $ws.Cells.Item($r, 2).Value = $val
}

# Update the view: scroll so row 458 is at the top, select A477 (first empty row)
$win = $excel.ActiveWindow
$win.ScrollRow = 458
$win.ScrollColumn = 1
$ws.Range("A477").Select()
